$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06446533333333333
$ws.Range("H2").Value = 0.193396
$ws.Range("I2").Value = 0.02693738696927793
$ws.Range("J2").Value = 0.02693738696927793
$ws.Range("M2").Value = 0.0002903333333333334
$ws.Range("N2").Value = 0.000871
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.00001871643511111111
$ws.Range("R2").Value = 0.000168447916
$ws.Range("S2").Value = 0.02693738696927793
$ws.Range("T2").Value = 0.02693738696927793

# --- Row 3 updates ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.843761666666667
$ws.Range("H3").Value = 5.531285
$ws.Range("I3").Value = 0.7704314695358874
$ws.Range("J3").Value = 0.7704314695358874
$ws.Range("M3").Value = 0.0002903333333333334
$ws.Range("N3").Value = 0.000871
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.0005353054705555556
$ws.Range("R3").Value = 0.004817749235000001
$ws.Range("S3").Value = 0.7704314695358874
$ws.Range("T3").Value = 0.7704314695358874

# --- Row 4 updates ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.4849276666666666
$ws.Range("H4").Value = 1.454783
$ws.Range("I4").Value = 0.2026311434948347
$ws.Range("J4").Value = 0.2026311434948347
$ws.Range("M4").Value = 0.0002903333333333334
$ws.Range("N4").Value = 0.000871
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.0001407906658888889
$ws.Range("R4").Value = 0.001267115993
$ws.Range("S4").Value = 0.2026311434948347
$ws.Range("T4").Value = 0.2026311434948347

# --- Remove old trailing rows 5-7 (data now only spans rows 1-4) ---
$ws.Rows("5:7").Delete()
